$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (BEATRIZ, 005046790): update Saldo value from 20104.34 to 29000
$ws.Range("C3").Value = 29000

# Delete rows 7-26 (MIRELLA through ALINE) first, from the bottom up,
# so the row numbers for the earlier deletion stay valid.
$ws.Range("A7:C26").EntireRow.Delete()

# Delete rows 4-5 (MERG 004214592/13545.88, THAYSA 004425261/8581.63)
$ws.Range("A4:C5").EntireRow.Delete()
